# "worked on monster properties"
# - Replace the "Student naam + Nummer" placeholder label with the
#   student's actual name/number.
# - Fill in the "Eigen Score" (self-assessment) column with the scores
#   the student is claiming for each criterion (matching the max. score
#   already listed in column D for the rows that are filled in).
# - Total the "Eigen Score" column with a SUM formula.
# - Leave the sheet scrolled/zoomed/selected where the student left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The label cell under the "Handtekening" (signature) block used to read
# "Student naam + Nummer" in both A40 and D40; replace it with the
# student's real name + number and drop the now-redundant A40 label.
$ws.Range("A40").ClearContents()
$ws.Range("D40").Value = "Simon Striekwold - 2137518"

# "1 - Applicatie Code" section: self-assessed ("Eigen Score") points,
# mirroring the max. score already entered in column D for each row.
$ws.Range("F17").Value = 5
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 3
$ws.Range("F20").Value = 6
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 5
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = 8

# "2 - Functionaliteit Applicatie" section: a few more self-assessed
# scores.
$ws.Range("F29").Value = 3
$ws.Range("F31").Value = 3
$ws.Range("F34").Value = 10

# Total of the "Eigen Score" column.
$ws.Range("F38").Formula = "=SUM(F17:F37)"

# Restore the view: zoomed to 85%, scrolled down to row 21, with F32
# selected (where the student was last working).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$win.Zoom = 85
$ws.Range("F32").Select()
